$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16 data, mirroring the style/format of the existing data rows (rows 2-15)
$rowIndex = 16

# Column A: index value, styled like the other index cells (bold/centered/bordered -> same style as A2:A15)
$ws.Cells.Item($rowIndex - 1, 1).Copy() | Out-Null
$ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item($rowIndex, 1).Value = 14

# Column B: label referencing the existing shared string "HexGrid-60degTilt5degRes"
$ws.Cells.Item($rowIndex, 2).Value = "HexGrid-60degTilt5degRes"

# Columns C through M: numeric averaged-intensity values
$ws.Cells.Item($rowIndex, 3).Value = 1.66636463905608
$ws.Cells.Item($rowIndex, 4).Value = 2.284038275547505
$ws.Cells.Item($rowIndex, 5).Value = 0.4213586778482891
$ws.Cells.Item($rowIndex, 6).Value = 1.66636463905608
$ws.Cells.Item($rowIndex, 7).Value = 1.319409511752389
$ws.Cells.Item($rowIndex, 8).Value = 0.7193734977715674
$ws.Cells.Item($rowIndex, 9).Value = 0.7084720708245341
$ws.Cells.Item($rowIndex, 10).Value = 2.284038275547505
$ws.Cells.Item($rowIndex, 11).Value = 1.352698476697897
$ws.Cells.Item($rowIndex, 12).Value = 1.509531557876988
$ws.Cells.Item($rowIndex, 13).Value = 1.186502778800061
